$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.684.43'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '3.782.56'
$ws.Range("E3").Value = '  -1.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '432.56'
$ws.Range("E5").Value = '  +5.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.54'
$ws.Range("E6").Value = '  +7.53%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.733'
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.153'
$ws.Range("E10").Value = '  -9.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000314'
$ws.Range("E11").Value = '  -14.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.84'
$ws.Range("E12").Value = '  +4.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.42'
$ws.Range("E13").Value = '  +4.28%  '
$ws.Range("D14").Value = '4.385.83'
$ws.Range("E14").Value = '  -1.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.96'
$ws.Range("E15").Value = '  -3.63%  '
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").Value = '3.816.76'
$ws.Range("E17").Value = '  -0.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.91'
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("E19").Value = '  +6.67%  '
$ws.Range("D20").Value = '66.798.26'
$ws.Range("E20").Value = '  -0.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '409.32'
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.74'
$ws.Range("E22").Value = '  +2.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.24'
$ws.Range("E23").Value = '  +6.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.48'
$ws.Range("E24").Value = '  +0.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '36.83'
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("E26").Value = '  +7.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.65'
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.63'
$ws.Range("E28").Value = '  +34.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.76'
$ws.Range("E29").Value = '  +3.90%  '
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '713.52'
$ws.Range("E30").Value = '  +4.53%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.83'
$ws.Range("E31").Value = '  +11.30%  '
$ws.Range("E32").Value = '  +10.54%  '
$ws.Range("E33").Value = '  +1.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '41.72'
$ws.Range("E34").Value = '  +8.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  +0.29%  '
$ws.Range("E37").Value = '  +26.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.06'
$ws.Range("E38").Value = '  +2.46%  '
$ws.Range("E39").Value = '  +3.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.74'
$ws.Range("E40").Value = '  +42.65%  '
$ws.Range("E41").Value = '  -5.19%  '
$ws.Range("B42").Value = 'ApeXProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.36'
$ws.Range("E42").Value = '  +7.30%  '
$ws.Range("E43").Value = '  +3.49%  '
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '0.0₃0674'
$ws.Range("E44").Value = '  -14.25%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.322'
$ws.Range("E46").Value = '  +9.38%  '
$ws.Range("E47").Value = '  +0.69%  '
$ws.Range("E48").Value = '  +5.40%  '
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '142.68'
$ws.Range("E50").Value = '  -4.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.82'
$ws.Range("E51").Value = '  +1.38%  '
